# fix: update test grid
#
# 1. Insert a new "Read Me" worksheet as the first sheet in the workbook,
#    containing a short note asking readers to report their usage of the
#    spreadsheet back to the author.
# 2. In "Areas Features Validations", tag a few Branching / Edit-Train
#    rows with their Test Name (column E) so they line up with the rest
#    of the EditAndBranching block.
# 3. Leave "Areas Features Validations" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. New "Read Me" sheet, moved to the front ------------------------
$readMe = $wb.Worksheets.Add()
$readMe.Name = "Read Me"
$readMe.Move($wb.Worksheets.Item(1))

$readMe.Range("A1").Value = "This spreadsheet was expensive to produce and to maintain…"
$readMe.Range("A2").Value = "therefore we should track its usefulness."
$readMe.Range("A3").Value = "Please send an email to Michael Skowronski (v-miskow) to let me know…"
$readMe.Range("A4").Value = "* if you use it"
$readMe.Range("A5").Value = "* how often you use it"
$readMe.Range("A6").Value = "* how you are using it"
$readMe.Range("A7").Value = "* if it gave you the information you wanted"
$readMe.Range("A9").Value = "How to use it…"

$readMe.Columns.Item(1).ColumnWidth = 66.6
$readMe.PageSetup.Orientation = 1

# --- 2. Fill in missing / updated Test Name values ----------------------
$grid = $wb.Worksheets.Item("Areas Features Validations")

$grid.Range("E71").Value = "EditAndBranching/Branching"
$grid.Range("E74").Value = "EditAndBranching/VerifyEditTrainingControlsAndLabels"
$grid.Range("E75").Value = "EditAndBranching/VerifyEditTrainingControlsAndLabels"
$grid.Range("E76").Value = "EditAndBranching/VerifyEditTrainingControlsAndLabels"

# --- 3. View state: "Areas Features Validations" stays the active tab --
$readMe.Activate()
$excel.ActiveWindow.Zoom = 130
$readMe.Range("A10").Select()

$grid.Activate()
$grid.Range("E80").Select()
